$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "269.80"
    "D3"  = "22.90"
    "D4"  = "6.389"
    "D5"  = "0.06211"
    "D6"  = "3.647"
    "D7"  = "6.734"
    "D8"  = "1.400"
    "D9"  = "0.8390"
    "D10" = "0.01361"
    "D11" = "0.1617"
    "D12" = "0.08258"
    "D13" = "0.03410"
    "D14" = "0.03208"
    "D15" = "0.09313"
    "D16" = "3.950"
    "D17" = "0.001737"
    "D18" = "0.04877"
    "D19" = "0.006300"
    "D20" = "0.005396"
    "D21" = "0.001100"
    "D22" = "0.0001515"
    "D23" = "3.758"
    "D24" = "2.366"
    "D25" = "0.3340"
    "D26" = "0.1248"
    "D27" = "0.0002708"
    "D40" = "0.04687"
    "D41" = "0.006932"
    "D42" = "0.1157"
    "D43" = "0.003451"
    "D44" = "0.01237"
    "D45" = "0.00006307"
    "D46" = "0.00000000757"
    "D47" = "0.7068"
    "D48" = "0.1306"
    "D49" = "0.00002120"
    "D50" = "0.01252"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text interpretation so trailing/leading zeros in the price
    # strings are preserved exactly (these are inline text values, not
    # numbers), then restore the original (default) cell style so only
    # the displayed text changes.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
